# Auto-generated edit script updating the cryptos price/volume table
# to reflect the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.034.73'
$ws.Range('E2').Value = '  +0.77%  '

# Row 3
$ws.Range('D3').Value = '3.420.64'
$ws.Range('E3').Value = '  -0.03%  '

# Row 4
$ws.Range('E4').Value = '  -0.21%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '406.26'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.37%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.80'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.54%  '

# Row 7
$ws.Range('E7').Value = '  -1.30%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.10%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.690'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.32%  '

# Row 10
$ws.Range('E10').Value = '  +3.40%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.91'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.10%  '

# Row 12
$ws.Range('E12').Value = '  -0.42%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.98'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.03%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.43'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.86%  '

# Row 15
$ws.Range('D15').Value = '3.412.76'
$ws.Range('E15').Value = '  -0.37%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '11.65'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.61%  '

# Row 17
$ws.Range('D17').Value = '61.855.72'
$ws.Range('E17').Value = '  +0.13%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.02'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.21%  '

# Row 19
$ws.Range('E19').Value = '  +8.92%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.18'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.16%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '84.09'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.87%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '313.90'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.73%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.81'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.56%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.16'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.67%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.75'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.80%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '29.66'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.26%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.99'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +6.42%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.14'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -5.54%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.77'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +8.25%  '

# Row 30
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.172'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.56%  '

# Row 31
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.115'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.27%  '

# Row 32
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '43.34'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.85%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.34'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.95%  '

# Row 34
$ws.Range('E34').Value = '  +0.10%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0488'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.71%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '51.51'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.30%  '

# Row 37
$ws.Range('E37').Value = '  -0.25%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.03'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.97%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.37'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.31%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.316'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +10.79%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '139.11'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.74%  '

# Row 42
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.98'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.46%  '

# Row 43
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.125'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.12%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.96'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.35%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.79'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.94%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.21'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.26%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.28'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.56%  '

# Row 48
$ws.Range('D48').Value = '2.105.42'
$ws.Range('E48').Value = '  -2.33%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.33'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.29%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.93'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.56%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.74'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +16.91%  '

